$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 348
$ws.Range("F3").Value = 3528
$ws.Range("F7").Value = 96
$ws.Range("F8").Value = 2184
$ws.Range("F10").Value = 190
$ws.Range("F12").Value = 1201
$ws.Range("F13").Value = 68
$ws.Range("F15").Value = 46
$ws.Range("F16").Value = 594
$ws.Range("F17").Value = 86
$ws.Range("F18").Value = 6247
$ws.Range("F20").Value = 7353
$ws.Range("F22").Value = 56387
$ws.Range("F23").Value = 4522
$ws.Range("F24").Value = 1043
$ws.Range("F25").Value = 882
$ws.Range("F26").Value = 439
$ws.Range("F27").Value = 95
$ws.Range("F29").Value = 600
$ws.Range("F30").Value = 3736
$ws.Range("F31").Value = 585
$ws.Range("F35").Value = 1237
$ws.Range("F36").Value = 1243
$ws.Range("F37").Value = 8
$ws.Range("F38").Value = 161
$ws.Range("F39").Value = 197
$ws.Range("F40").Value = 1072
$ws.Range("F41").Value = 710
$ws.Range("F43").Value = 772
$ws.Range("F44").Value = 180
$ws.Range("F46").Value = 172
$ws.Range("F47").Value = 7
$ws.Range("F49").Value = 2475

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value = 176
$ws.Range("F11").Value = 51
$ws.Range("F12").Value = 116
$ws.Range("F16").Value = 7495
$ws.Range("F38").Value = 36
$ws.Range("F43").Value = 94

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 2299
$ws.Range("F5").Value = 1557
$ws.Range("F7").Value = 665
$ws.Range("F8").Value = 2350
$ws.Range("F9").Value = 9344
$ws.Range("F10").Value = 1686
$ws.Range("F11").Value = 161
$ws.Range("F15").Value = 186
$ws.Range("F16").Value = 329

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 3528
$ws.Range("F3").Value = 8243
$ws.Range("F4").Value = 1557
$ws.Range("F5").Value = 2350
$ws.Range("F7").Value = 1686
$ws.Range("F8").Value = 161
$ws.Range("F11").Value = 96
$ws.Range("F16").Value = 190
$ws.Range("F17").Value = 68
$ws.Range("F18").Value = 46
$ws.Range("F19").Value = 594
$ws.Range("F20").Value = 86
$ws.Range("F21").Value = 56387
$ws.Range("F23").Value = 4522
$ws.Range("F24").Value = 51
$ws.Range("F25").Value = 439
$ws.Range("F26").Value = 600
$ws.Range("F27").Value = 116
$ws.Range("F29").Value = 3736
$ws.Range("F30").Value = 585
$ws.Range("F35").Value = 1237
$ws.Range("F37").Value = 329
$ws.Range("F39").Value = 161
$ws.Range("F40").Value = 1072
$ws.Range("F41").Value = 710
$ws.Range("F42").Value = 772
$ws.Range("F43").Value = 180
$ws.Range("F44").Value = 172
$ws.Range("F45").Value = 7
$ws.Range("F46").Value = 36
$ws.Range("F49").Value = 2475
